$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Range("K2").Value = 5020
$ws.Range("J3").Value = 8077
$ws.Range("K3").Value = 5146
$ws.Range("K4").Value = 1069
$ws.Range("K5").Value = 364
$ws.Range("K6").Value = 5787
$ws.Range("J7").Value = 29298
$ws.Range("K7").Value = 17386

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Range("K2").Value = 56
$ws.Range("K3").Value = 55
$ws.Range("K7").Value = 225

$ws = $wb.Worksheets.Item("Austin")
$ws.Range("K2").Value = 324
$ws.Range("K3").Value = 345
$ws.Range("K6").Value = 395
$ws.Range("K7").Value = 1163

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Range("K3").Value = 141
$ws.Range("K6").Value = 89
$ws.Range("K7").Value = 384

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Range("K2").Value = 201
$ws.Range("K6").Value = 214
$ws.Range("K7").Value = 735

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Range("K3").Value = 105
$ws.Range("K7").Value = 297

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Range("K2").Value = 169
$ws.Range("J3").Value = 301
$ws.Range("J7").Value = 901
$ws.Range("K7").Value = 587

$ws = $wb.Worksheets.Item("New City")
$ws.Range("K2").Value = 123
$ws.Range("K3").Value = 100
$ws.Range("K6").Value = 155
$ws.Range("K7").Value = 402

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Range("K2").Value = 150
$ws.Range("K4").Value = 65
$ws.Range("K6").Value = 131
$ws.Range("K7").Value = 519
$ws.Range("K8").Value = 1163
$ws.Range("K15").Value = 175
$ws.Range("K18").Value = 116
$ws.Range("K19").Value = 520
$ws.Range("K27").Value = 160
$ws.Range("K29").Value = 935
$ws.Range("K33").Value = 735
$ws.Range("K34").Value = 93
$ws.Range("K36").Value = 229
$ws.Range("J37").Value = 901
$ws.Range("K37").Value = 587
$ws.Range("K39").Value = 22
$ws.Range("K43").Value = 153
$ws.Range("K50").Value = 85
$ws.Range("J51").Value = 369
$ws.Range("K51").Value = 222
$ws.Range("K52").Value = 450
$ws.Range("K53").Value = 225
$ws.Range("K54").Value = 341
$ws.Range("K55").Value = 198
$ws.Range("K57").Value = 64
$ws.Range("J63").Value = 112
$ws.Range("K63").Value = 48
$ws.Range("K65").Value = 402
$ws.Range("K67").Value = 668
$ws.Range("K72").Value = 84
$ws.Range("K73").Value = 150
$ws.Range("K76").Value = 240
$ws.Range("K79").Value = 425
$ws.Range("K83").Value = 384
$ws.Range("K85").Value = 805
$ws.Range("K86").Value = 118
$ws.Range("K87").Value = 29
$ws.Range("K89").Value = 251
$ws.Range("K91").Value = 186
$ws.Range("K92").Value = 66
$ws.Range("K93").Value = 66
$ws.Range("K95").Value = 297
$ws.Range("K96").Value = 185
$ws.Range("K97").Value = 138
$ws.Range("K100").Value = 33
$ws.Range("J101").Value = 29298
$ws.Range("K101").Value = 17386

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Range("K2").Value = 193
$ws.Range("K3").Value = 233
$ws.Range("K7").Value = 668

$ws = $wb.Worksheets.Item("Loop")
$ws.Range("K3").Value = 89
$ws.Range("K7").Value = 341

$ws = $wb.Worksheets.Item("Englewood")
$ws.Range("K3").Value = 332
$ws.Range("K7").Value = 935

$ws = $wb.Worksheets.Item("Chatham")
$ws.Range("K5").Value = 17
$ws.Range("K7").Value = 520

$ws = $wb.Worksheets.Item("River North")
$ws.Range("K2").Value = 50
$ws.Range("K3").Value = 44
$ws.Range("K7").Value = 240

$ws = $wb.Worksheets.Item("Ashburn")
$ws.Range("K6").Value = 38
$ws.Range("K7").Value = 131

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Range("K2").Value = 169
$ws.Range("K3").Value = 200

$ws = $wb.Worksheets.Item("Lower West Side")
$ws.Range("K3").Value = 58
$ws.Range("K6").Value = 69
$ws.Range("K7").Value = 198

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Range("K2").Value = 58
$ws.Range("K7").Value = 185

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Range("K2").Value = 48
$ws.Range("K7").Value = 186

$ws = $wb.Worksheets.Item("Roseland")
$ws.Range("K6").Value = 105
$ws.Range("K7").Value = 425

$ws = $wb.Worksheets.Item("Calumet Heights")
$ws.Range("K3").Value = 36
$ws.Range("K7").Value = 116

$ws = $wb.Worksheets.Item("Grand Boulevard")
$ws.Range("K2").Value = 92
$ws.Range("K7").Value = 229

$ws = $wb.Worksheets.Item("West Lawn")
$ws.Range("K3").Value = 14
$ws.Range("K7").Value = 66

$ws = $wb.Worksheets.Item("Wrigleyville")
$ws.Range("K2").Value = 6
$ws.Range("K7").Value = 33

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Range("K3").Value = 169
$ws.Range("K7").Value = 519

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Range("K2").Value = 33
$ws.Range("K7").Value = 93

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Range("K3").Value = 43
$ws.Range("K7").Value = 175

$ws = $wb.Worksheets.Item("Lincoln Square")
$ws.Range("K6").Value = 43
$ws.Range("K7").Value = 85

$ws = $wb.Worksheets.Item("Greektown")
$ws.Range("K5").Value = 13
$ws.Range("K6").Value = 22

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Range("K3").Value = 38
$ws.Range("K7").Value = 150

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Range("K2").Value = 49
$ws.Range("K3").Value = 35
$ws.Range("K4").Value = 12
$ws.Range("K7").Value = 150

$ws = $wb.Worksheets.Item("West Town")
$ws.Range("K6").Value = 83
$ws.Range("K7").Value = 138

$ws = $wb.Worksheets.Item("West Elsdon")
$ws.Range("K2").Value = 20
$ws.Range("K7").Value = 66

$ws = $wb.Worksheets.Item("Uptown")
$ws.Range("K3").Value = 78
$ws.Range("K7").Value = 251

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Range("K2").Value = 43
$ws.Range("K7").Value = 160

$ws = $wb.Worksheets.Item("Streeterville")
$ws.Range("K4").Value = 49
$ws.Range("K7").Value = 118

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Range("K2").Value = 65
$ws.Range("J4").Value = 35
$ws.Range("J7").Value = 369
$ws.Range("K7").Value = 222

$ws = $wb.Worksheets.Item("Mckinley Park")
$ws.Range("K6").Value = 31
$ws.Range("K7").Value = 64

$ws = $wb.Worksheets.Item("Hyde Park")
$ws.Range("K2").Value = 30
$ws.Range("K7").Value = 153

$ws = $wb.Worksheets.Item("South Shore")
$ws.Range("K3").Value = 269
$ws.Range("K4").Value = 49
$ws.Range("K6").Value = 194
$ws.Range("K7").Value = 805

$ws = $wb.Worksheets.Item("Old Town")
$ws.Range("K6").Value = 43
$ws.Range("K7").Value = 84

$ws = $wb.Worksheets.Item("Little Village")
$ws.Range("K6").Value = 167
$ws.Range("K7").Value = 450

$ws = $wb.Worksheets.Item("Archer Heights")
$ws.Range("K6").Value = 25
$ws.Range("K7").Value = 65

$ws = $wb.Worksheets.Item("Ukrainian Village")
$ws.Range("K6").Value = 16
$ws.Range("K7").Value = 29
